# Update Name of Algo
# Apply updated KNN-imputed values to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value  = -21.722
$ws.Range("D3").Value  = -7.824000000000001

$ws.Range("E19").Value = 16.441

$ws.Range("A21").Value = -20.469

$ws.Range("A23").Value = -20.428
$ws.Range("D24").Value = -7.456999999999999
$ws.Range("E24").Value = 17.302

$ws.Range("A25").Value = -21.664

$ws.Range("B27").Value = 5.723000000000001

$ws.Range("E30").Value = 16.26

$ws.Range("B31").Value = 6.106
$ws.Range("E31").Value = 16.246

$ws.Range("E33").Value = 17.229

$ws.Range("B39").Value = 8.132000000000001

$ws.Range("B48").Value = 5.323

$ws.Range("B51").Value = 6.065

$ws.Range("B52").Value = 5.398000000000001

$ws.Range("A53").Value = -21.915

$ws.Range("B55").Value = 4.697
$ws.Range("E55").Value = 16.354

$ws.Range("B56").Value = 5.003

$ws.Range("A57").Value = -21.303
$ws.Range("B57").Value = 6.431999999999999
$ws.Range("D57").Value = -8.231999999999999

$ws.Range("A59").Value = -22.266

$ws.Range("D61").Value = -7.704000000000001

$ws.Range("E65").Value = 17.336

$ws.Range("A69").Value = -21.649

$ws.Range("D70").Value = -6.857000000000001
$ws.Range("E70").Value = 17.746

$ws.Range("B73").Value = 6.751

$ws.Range("E75").Value = 16.59

$ws.Range("A79").Value = -21.137

$ws.Range("A83").Value = -22
$ws.Range("E83").Value = 16.826

$ws.Range("D86").Value = -8.218

$ws.Range("B89").Value = 5.986999999999999

$ws.Range("B90").Value = 5.671

$ws.Range("A93").Value = -21.498

$ws.Range("E96").Value = 16.454

$ws.Range("E97").Value = 16.767

$ws.Range("D98").Value = -7.973999999999999

$ws.Range("D100").Value = -8.022000000000002

$ws.Range("D102").Value = -8.003000000000002
